$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 188-189, pushing the existing rows
# (old 188 onward) down by two positions (old 188 -> new 190, etc.)
$ws.Rows("188:189").Insert()

# Populate the two newly inserted rows with the new weekly data
# Row 188: "1a amarillo"
$ws.Range("A188").Value = 11
$ws.Range("B188").Value = "Vega Monumental Concepción"
$ws.Range("C188").Value = "Bíobío"
$ws.Range("D188").Value = "2021-09-29"
$ws.Range("E188").Value = 8
$ws.Range("F188").Value = "Fruta"
$ws.Range("G188").Value = 100102
$ws.Range("H188").Value = "Cítricos"
$ws.Range("I188").Value = 100102003
$ws.Range("J188").Value = "Limón"
$ws.Range("K188").Value = "Sin especificar"
$ws.Range("L188").Value = "1a amarillo"
$ws.Range("M188").Value = 600
$ws.Range("N188").Value = 6000
$ws.Range("O188").Value = 6500
$ws.Range("P188").Value = 6250
$ws.Range("Q188").Value = "$/malla 16 kilos"
$ws.Range("R188").Value = "Provincia de Melipilla"
$ws.Range("S188").Value = 391
$ws.Range("T188").Value = 16

# Row 189: "2a amarillo"
$ws.Range("A189").Value = 11
$ws.Range("B189").Value = "Vega Monumental Concepción"
$ws.Range("C189").Value = "Bíobío"
$ws.Range("D189").Value = "2021-09-29"
$ws.Range("E189").Value = 8
$ws.Range("F189").Value = "Fruta"
$ws.Range("G189").Value = 100102
$ws.Range("H189").Value = "Cítricos"
$ws.Range("I189").Value = 100102003
$ws.Range("J189").Value = "Limón"
$ws.Range("K189").Value = "Sin especificar"
$ws.Range("L189").Value = "2a amarillo"
$ws.Range("M189").Value = 300
$ws.Range("N189").Value = 5500
$ws.Range("O189").Value = 5500
$ws.Range("P189").Value = 5500
$ws.Range("Q189").Value = "$/malla 16 kilos"
$ws.Range("R189").Value = "Provincia de Melipilla"
$ws.Range("S189").Value = 344
$ws.Range("T189").Value = 16

# Make sure the D column keeps the date number-format used
# elsewhere in that column (style index 2 in the original workbook).
$ws.Range("D188:D189").NumberFormat = $ws.Range("D190").NumberFormat
